# Roll the reporting window of the income-statement workbook forward by one
# fiscal year (drop FY1396/12, add FY1401/12) and reset every financial
# figure to 0 (or the "-" placeholder) per the new read_price algorithm.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: fiscal-period headers (columns D..H) ---
$ws.Cells.Item(8, 4).Value = "12 ماهه منتهی به 1397/12"
$ws.Cells.Item(8, 5).Value = "12 ماهه منتهی به 1398/12"
$ws.Cells.Item(8, 6).Value = "12 ماهه منتهی به 1399/12"
$ws.Cells.Item(8, 7).Value = "12 ماهه منتهی به 1400/12"
$ws.Cells.Item(8, 8).Value = "12 ماهه منتهی به 1401/12"

# --- Row 9: publish dates (columns D..H) ---
$ws.Cells.Item(9, 4).Value = "1399-04-21 (13)"
$ws.Cells.Item(9, 5).Value = "1400-02-30 (7)"
$ws.Cells.Item(9, 6).Value = "1401-04-18 (7)"
$ws.Cells.Item(9, 7).Value = "1402-02-23 (8)"

# H9 is the bare "1402-02-23" (no trailing "(n)" qualifier), which Excel's
# COM value coercion would otherwise auto-parse as a date literal and store
# as a serial number. Route it through a text formula, then flatten the
# formula down to a literal value in place via copy / paste-special so the
# stored cell keeps its plain-text (shared-string) representation and its
# original style.
$ws.Cells.Item(9, 8).Formula = "=""1402-02-23"""
$ws.Cells.Item(9, 8).Copy()
$ws.Cells.Item(9, 8).PasteSpecial(-4163)

# --- Rows 11..27: reset financial data to 0, preserving the "-" rows ---
$dashRows = @(15, 23)
for ($r = 11; $r -le 27; $r++) {
    for ($c = 4; $c -le 8; $c++) {
        if ($dashRows -contains $r) {
            $ws.Cells.Item($r, $c).Value = "-"
        } else {
            $ws.Cells.Item($r, $c).Value = 0
        }
    }
}
